$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Cadastrado" (column B) values
$ws.Range("B2").Value = 102012.86
$ws.Range("B3").Value = 906814.4
$ws.Range("B4").Value = 1768046.15
$ws.Range("B5").Value = 2869200.58
$ws.Range("B6").Value = 4501569.72
$ws.Range("B7").Value = 1906073.08

# Add new "Sem Cadastro" column (C) with header + values
$ws.Range("C1").Value = "Sem Cadastro"
$ws.Range("C2").Value = 13664.65
$ws.Range("C3").Value = 46293.98
$ws.Range("C4").Value = 18659.04
$ws.Range("C5").Value = 16773.44
$ws.Range("C6").Value = 15863.05
$ws.Range("C7").Value = 4029.02

# Copy the header formatting from B1 onto C1 (matches existing style, no new style entries)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
